$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = 27556940.46
$ws.Range("P2").Value = 45.3456806125
$ws.Range("Q2").Value = 105540785.91
$ws.Range("R2").Value = 173.6701785314
$ws.Range("S2").Value = 26534108.61
$ws.Range("T2").Value = 43.6625835191
$ws.Range("U2").Value = 49467506.93
$ws.Range("V2").Value = 81.4001022066
$ws.Range("Y2").Value = 2064921.99
$ws.Range("Z2").Value = 3.3978842167
$ws.Range("AA2").Value = -16241667
$ws.Range("AB2").Value = -26.7260962974
$ws.Range("AC2").Value = 60770816.73
$ws.Range("AD2").Value = ""
